$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9588997364044189
$ws.Range("B1").Value = 2.213033676147461
$ws.Range("C1").Value = 8.332962989807129
$ws.Range("D1").Value = 1.702846765518188
$ws.Range("E1").Value = 1.337262630462646
